$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Set cells in the order that introduces new shared strings matching
# the target string-table layout: B23, B24, C23, D23, then the rest.
$ws.Range("B23").Value = "Check marriage if before 14"
$ws.Range("B24").Value = "check if marry twice at same time"
$ws.Range("C23").Value = "WYK"
$ws.Range("D23").Value = "coding"

$ws.Range("A23").Value = "US10"
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 50

$ws.Range("A24").Value = "US11"
$ws.Range("C24").Value = "WYK"
$ws.Range("D24").Value = "coding"
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 50

[void]$ws.Range("B21").Select()
